$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New Cypher query text (shared strings 10-13 in the final workbook)
# ---------------------------------------------------------------------------

$statQuery = @"
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Targeted-Capture" IN es
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s:study)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Targeted-Capture" IN es
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Targeted-Capture" IN es
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS ``Files``
"@

$participantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (samp)<--(f:file)
WITH p, samp, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Targeted-Capture" IN es
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

$samplesQuery = @"
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Targeted-Capture" IN es
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(p.participant_id,'') as ``Participant ID``,
    coalesce(s.study_name, '') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(samp.sample_tumor_status,'') as ``Tumor``,
    coalesce(samp.sample_type,'') as ``Analyte Type``
ORDER BY samp.sample_id LIMIT 100
"@

$filesQuery = @"
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Targeted-Capture" IN es
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as ``File Name``,
    coalesce(s.study_name,'') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(p.participant_id, '') as ``Participant ID``,
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(f.file_type, '') as ``File Type``
ORDER BY f.file_name LIMIT 100
"@

# ---------------------------------------------------------------------------
# Update the cell contents for the three tabs (Participants, Samples, Files)
# Note: B/C columns swap meaning on row 2 vs rows 3/4 in original layout;
# all rows now show <tab query> in column B and the shared StatQuery in C.
# Write order matters for shared-string de-dup ordering, so write the
# StatQuery text first (it ends up first in the shared-string table).
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = $statQuery
$ws.Range("B2").Value = $participantsQuery

$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery

$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# ---------------------------------------------------------------------------
# Formatting: bump the query-column font size from 12 to 14, keep wrap text.
# ---------------------------------------------------------------------------
$ws.Range("B2:C4").Font.Size = 14
$ws.Range("B2:C4").WrapText = $true

$ws.Range("B5:C5").Font.Size = 14
$ws.Range("B5:C5").WrapText = $true
$ws.Range("C6").Font.Size = 14
$ws.Range("C6").WrapText = $true

# Row heights grow to the max (409.5) to fit the longer wrapped queries.
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# ---------------------------------------------------------------------------
# Selection / view: scrolled down one row, new active cell B4.
# ---------------------------------------------------------------------------
$ws.Range("A3").Select()
$ws.Range("B4").Select()
